$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("known issues")

# New "known issue" row (row 7): unit-conversion issue between parent/child layers.
$ws.Range("A7").Value = "need to enforce unit conversion between parent and child layers (and log each instance of unit conversion)"
$ws.Range("B7").Value = "Currently unit difference is not checked between parent and child layers (e.g., parent layer uses 'MJ' for electricity input, but electricity generation (child layer) is in the unit of 'kWh')"
$ws.Range("C7").Value = "The amount of child layer activity may not be correct (as current code multiplies the amt of parent activity with the amt of child activity, without considering the difference in units)"

# Row grew tall to fit the wrapped text (matches the author's manual entry in Excel).
$ws.Rows.Item(7).RowHeight = 68

# Stray character left in J12 (cursor landed there), also becomes the new selection.
$ws.Range("J12").Value = "ß"
$ws.Range("J12").Select()
